# feat(KeyedList): support enum as key type
#
# The "Loader" sheet's config-name column (ServerConfName / [Conf]<string>)
# is turned into a config-type enum column (ServerConfType /
# [Conf]<enum<.ConfType>>), and the per-server config values switch from
# the removed bespoke config-sheet names (HeadFrameConf/ActivityConf/
# ChapterConf/CollectionConf/ExchangeConf/SectionConf/MatchConf) to the
# new CONF_TYPE_* enum members. The "@TABLEAU" / "Sheet2" meta sheets
# also relabel their "Loader" row to the generic "进程名" field name.

$wb = $excel.ActiveWorkbook

# --- Loader sheet: rename header + retarget values to the new enum type ---
$loader = $wb.Worksheets.Item("Loader")

$loader.Range("B1").Value = "ServerConfType"
$loader.Range("B2").Value = "[Conf]<enum<.ConfType>>"

$loader.Range("B5").Value = "CONF_TYPE_CLOUD"
$loader.Range("B7").Value = "CONF_TYPE_CLOUD"
$loader.Range("B8").Value = "CONF_TYPE_CLOUD"
$loader.Range("B9").Value = "CONF_TYPE_LOCAL"
$loader.Range("B10").Value = "CONF_TYPE_LOCAL"
$loader.Range("B11").Value = "CONF_TYPE_LOCAL"
$loader.Range("B12").Value = "CONF_TYPE_REMOTE"
$loader.Range("B13").Value = "CONF_TYPE_UNKNOWN"

# Column B got wider to fit the longer enum-type text (20.625 -> 25.875
# character-widths; 25 is what the COM ColumnWidth rounding needs to land
# back on that stored width).
$loader.Columns.Item(2).ColumnWidth = 25

# Active selection moved up one row (D16 -> D15).
[void]$loader.Range("D15").Select()

# --- Meta sheets: the "Loader" label becomes the generic "进程名" field ---
$tableau = $wb.Worksheets.Item("@TABLEAU")
$tableau.Range("A5").Value = "进程名"

$sheet2 = $wb.Worksheets.Item("Sheet2")
$sheet2.Range("A5").Value = "进程名"
